# Applies the cryptos-list refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.335.08"
$ws.Range("E2").Value = "  -2.04%  "

$ws.Range("D3").Value = "'2.515.29"
$ws.Range("E3").Value = "  -3.58%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'558.54"
$ws.Range("E5").Value = "  -2.66%  "

$ws.Range("D6").Value = "'148.39"
$ws.Range("E6").Value = "  -4.19%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "'0.615"
$ws.Range("E8").Value = "  -0.86%  "

$ws.Range("D9").Value = "'2.514.82"
$ws.Range("E9").Value = "  -3.53%  "

$ws.Range("D10").Value = "'0.110"
$ws.Range("E10").Value = "  -6.20%  "

$ws.Range("D11").Value = "'5.47"
$ws.Range("E11").Value = "  -5.95%  "

$ws.Range("E12").Value = "  -0.90%  "

$ws.Range("D13").Value = "'0.365"
$ws.Range("E13").Value = "  -4.32%  "

$ws.Range("D14").Value = "'26.55"
$ws.Range("E14").Value = "  -5.96%  "

$ws.Range("D15").Value = "'2.973.34"
$ws.Range("E15").Value = "  -3.39%  "

$ws.Range("D16").Value = "'0.0000169"
$ws.Range("E16").Value = "  -5.33%  "

$ws.Range("D17").Value = "'62.250.08"
$ws.Range("E17").Value = "  -1.91%  "

$ws.Range("D18").Value = "'2.535.21"
$ws.Range("E18").Value = "  -2.71%  "

$ws.Range("D19").Value = "'11.35"
$ws.Range("E19").Value = "  -5.21%  "

$ws.Range("D20").Value = "'7.11"
$ws.Range("E20").Value = "  -5.13%  "

$ws.Range("D21").Value = "'4.27"
$ws.Range("E21").Value = "  -5.89%  "

$ws.Range("D22").Value = "'328.82"
$ws.Range("E22").Value = "  -4.13%  "

$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").Value = "'65.25"
$ws.Range("E24").Value = "  -2.68%  "

$ws.Range("D25").Value = "'1.77"
$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("D26").Value = "'0.0000106"
$ws.Range("E26").Value = "  -1.19%  "

$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "'2.648.67"
$ws.Range("E27").Value = "  -2.95%  "

$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").Value = "'1.54"
$ws.Range("E28").Value = "  -0.94%  "

$ws.Range("D29").Value = "'8.66"
$ws.Range("E29").Value = "  -5.29%  "

$ws.Range("D30").Value = "'549.88"
$ws.Range("E30").Value = "  -6.23%  "

$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'7.86"
$ws.Range("E32").Value = "  -0.63%  "

$ws.Range("E33").Value = "  -3.37%  "

$ws.Range("D34").Value = "'1.94"
$ws.Range("E34").Value = "  -5.66%  "

$ws.Range("D35").Value = "'1.61"
$ws.Range("E35").Value = "  -7.13%  "

$ws.Range("D36").Value = "'6.07"
$ws.Range("E36").Value = "  -7.12%  "

$ws.Range("D37").Value = "'4.93"
$ws.Range("E37").Value = "  -8.47%  "

$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("D39").Value = "'0.384"
$ws.Range("E39").Value = "  -4.84%  "

$ws.Range("D40").Value = "'18.87"
$ws.Range("E40").Value = "  -3.91%  "

$ws.Range("D41").Value = "'149.25"
$ws.Range("E41").Value = "  -3.16%  "

$ws.Range("D42").Value = "'1.73"
$ws.Range("E42").Value = "  -7.14%  "

$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("D44").Value = "'40.99"
$ws.Range("E44").Value = "  -1.10%  "

$ws.Range("D45").Value = "'2.38"
$ws.Range("E45").Value = "  -2.55%  "

$ws.Range("D46").Value = "'150.68"
$ws.Range("E46").Value = "  -3.08%  "

$ws.Range("D47").Value = "'3.69"
$ws.Range("E47").Value = "  -5.11%  "

$ws.Range("D48").Value = "'21.73"
$ws.Range("E48").Value = "  -6.44%  "

$ws.Range("D49").Value = "'0.0549"
$ws.Range("E49").Value = "  -6.46%  "

$ws.Range("D50").Value = "'0.599"
$ws.Range("E50").Value = "  -4.43%  "

$ws.Range("D51").Value = "'0.0963"
$ws.Range("E51").Value = "  -3.70%  "
